# "expanded years to 2029"
# - Update the End Year value on the "Coupling Parameters" sheet from 2025 to 2029.
# - Make "Coupling Parameters" the active/selected sheet (was "Import Priorities").
# - Update the selection on "Coupling Parameters" to D9.

$wb = $excel.ActiveWorkbook

$couplingSheet = $wb.Worksheets.Item("Coupling Parameters")

# End Year (row 3, column B) grows from 2025 to 2029.
$couplingSheet.Range("B3").Value = 2029

# Switch the active sheet to "Coupling Parameters" and move the selection to D9.
$couplingSheet.Activate()
$couplingSheet.Range("D9").Select()
